$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6855
$ws1.Range("F4").Value = 31
$ws1.Range("F5").Value = 446
$ws1.Range("F7").Value = 6551
$ws1.Range("F9").Value = 199
$ws1.Range("F10").Value = 1290
$ws1.Range("F13").Value = 398
$ws1.Range("F14").Value = 135
$ws1.Range("F15").Value = 18
$ws1.Range("F16").Value = 383
$ws1.Range("F18").Value = 9
$ws1.Range("F19").Value = 4906
$ws1.Range("F20").Value = 93
$ws1.Range("F21").Value = 93
$ws1.Range("F22").Value = 360
$ws1.Range("F23").Value = 207
$ws1.Range("F24").Value = 162

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 45

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6855
$ws4.Range("F3").Value = 89
$ws4.Range("F4").Value = 31
$ws4.Range("F5").Value = 446
$ws4.Range("F6").Value = 149
$ws4.Range("F7").Value = 6551
$ws4.Range("F8").Value = 58
$ws4.Range("F9").Value = 199
$ws4.Range("F10").Value = 1290
$ws4.Range("F13").Value = 398
$ws4.Range("F14").Value = 135
$ws4.Range("F16").Value = 383
$ws4.Range("F18").Value = 9
$ws4.Range("F19").Value = 4906
$ws4.Range("F20").Value = 45
$ws4.Range("F21").Value = 93
$ws4.Range("F22").Value = 93
$ws4.Range("F23").Value = 360
$ws4.Range("F24").Value = 207
$ws4.Range("F25").Value = 162
